# Fruta / hortaliza, semanal
# Insert 14 new weekly price rows for "Vega Modelo de Temuco" - Manzana,
# dated 2021-09-09 (serial 44448), just above the existing row 1170 block
# (which shifts down to become rows 1184:1202).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing rows 1170:1188 down by 14 rows, opening up a gap at
# 1170:1183 for the new records.
$ws.Rows("1170:1183").Insert()

# Columns that are constant across this market/product block.
$colA = 10
$colB = 'Vega Modelo de Temuco'
$colC = 'La Araucanía'
$colE = 9
$colF = 'Fruta'
$colG = 100104
$colH = 'Frutos de pepita'
$colI = 100104002
$colJ = 'Manzana'

# New rows: D=fecha(serial), K=variedad, L=calidad, M=volumen, N=precio min,
# O=precio max, P=precio promedio ponderado, Q=unidad, R=origen,
# S=precio $/kg, T=kg/unidad
$newRows = @(
  @(44448, 'Fuji royal', 'Calibre 80', 450, 18000, 19000, 18556, '$/caja 18 kilos embalada', 'Región del Maule', 1031, 18),
  @(44448, 'Fuji royal', 'Primera', 330, 11000, 12000, 11455, '$/bandeja 15 kilos granel', 'Región del Maule', 764, 15),
  @(44448, 'Fuji royal', 'Primera', 10, 220000, 220000, 220000, '$/bins (400 kilos)', 'Región del Maule', 550, 400),
  @(44448, 'Fuji royal', 'Segunda', 5, 180000, 180000, 180000, '$/bins (400 kilos)', 'Región del Maule', 450, 400),
  @(44448, 'Granny Smith', 'Calibre 120', 140, 16000, 16000, 16000, '$/caja 18 kilos embalada', 'Región del Maule', 889, 18),
  @(44448, 'Granny Smith', 'Calibre 80', 300, 18000, 19000, 18667, '$/caja 18 kilos embalada', 'Región del Maule', 1037, 18),
  @(44448, 'Granny Smith', 'Primera', 120, 12000, 12000, 12000, '$/bandeja 15 kilos granel', 'Región del Maule', 800, 15),
  @(44448, 'Granny Smith', 'Segunda', 50, 8000, 8000, 8000, '$/bandeja 15 kilos granel', 'Región del Maule', 533, 15),
  @(44448, 'Granny Smith', 'Segunda', 5, 180000, 180000, 180000, '$/bins (400 kilos)', 'Región del Maule', 450, 400),
  @(44448, 'Pink Lady', 'Primera', 3, 160000, 160000, 160000, '$/bins (400 kilos)', 'Región del Maule', 400, 400),
  @(44448, 'Richared Delicious', 'Calibre 80', 220, 19000, 19000, 19000, '$/caja 18 kilos embalada', 'Región del Maule', 1056, 18),
  @(44448, 'Royal Gala', 'Especial', 3, 250000, 250000, 250000, '$/bins (400 kilos)', 'Región del Maule', 625, 400),
  @(44448, 'Royal Gala', 'Primera', 100, 12000, 12000, 12000, '$/bandeja 15 kilos granel', 'Región del Maule', 800, 15),
  @(44448, 'Royal Gala', 'Primera', 5, 150000, 150000, 150000, '$/bins (400 kilos)', 'Región del Maule', 375, 400)
)

$r = 1170
foreach ($row in $newRows) {
  $ws.Cells.Item($r, 1).Value = $colA
  $ws.Cells.Item($r, 2).Value = $colB
  $ws.Cells.Item($r, 3).Value = $colC
  $ws.Cells.Item($r, 4).Value = $row[0]
  $ws.Cells.Item($r, 5).Value = $colE
  $ws.Cells.Item($r, 6).Value = $colF
  $ws.Cells.Item($r, 7).Value = $colG
  $ws.Cells.Item($r, 8).Value = $colH
  $ws.Cells.Item($r, 9).Value = $colI
  $ws.Cells.Item($r, 10).Value = $colJ
  $ws.Cells.Item($r, 11).Value = $row[1]
  $ws.Cells.Item($r, 12).Value = $row[2]
  $ws.Cells.Item($r, 13).Value = $row[3]
  $ws.Cells.Item($r, 14).Value = $row[4]
  $ws.Cells.Item($r, 15).Value = $row[5]
  $ws.Cells.Item($r, 16).Value = $row[6]
  $ws.Cells.Item($r, 17).Value = $row[7]
  $ws.Cells.Item($r, 18).Value = $row[8]
  $ws.Cells.Item($r, 19).Value = $row[9]
  $ws.Cells.Item($r, 20).Value = $row[10]
  $r = $r + 1
}
